$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend header row 1 with new columns P1, Q1, copying the existing bold/
# border/center formatting already used across the rest of row 1 (style
# index 1) from the neighboring O1 cell.
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# Flip values in columns I, K, M, O for data rows 2-25, and append the new
# P/Q data columns (value 2, unstyled like the rest of the data columns).
For ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value = 2   # I -> 2
    $ws.Cells.Item($r, 11).Value = 1  # K -> 1
    $ws.Cells.Item($r, 13).Value = 2  # M -> 2
    $ws.Cells.Item($r, 15).Value = 1  # O -> 1

    $ws.Cells.Item($r, 16).Value = 2  # P
    $ws.Cells.Item($r, 17).Value = 2  # Q
}
